$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 42; this shifts the existing rows 42..134
# down to 43..135 (matching the new dimension A1:R135 from the diff).
$ws.Rows(42).Insert()

# Populate the newly inserted row 42 with the new price record that the
# diff introduces (a weekly "Espinaca" quote for Vega Modelo de Temuco).
$ws.Cells.Item(42, 1).Value = 10
$ws.Cells.Item(42, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(42, 3).Value = "La Araucanía"
$ws.Cells.Item(42, 4).Value = 44662
$ws.Cells.Item(42, 5).Value = 9
$ws.Cells.Item(42, 6).Value = 100112012
$ws.Cells.Item(42, 7).Value = "Espinaca"
$ws.Cells.Item(42, 8).Value = "Sin especificar"
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 40
$ws.Cells.Item(42, 11).Value = 9000
$ws.Cells.Item(42, 12).Value = 9000
$ws.Cells.Item(42, 13).Value = 9000
$ws.Cells.Item(42, 14).Value = "`$/docena de atados"
$ws.Cells.Item(42, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(42, 16).Value = 3000
$ws.Cells.Item(42, 17).Value = 3
$ws.Cells.Item(42, 18).Value = "Hortaliza"
